$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Update the Cypher query text in cell B4: the WHERE clause filter conditions
# changed from filtering on diagnosis stage/file_format('rtf') to filtering on
# demo.sex/file_type('Pathology Report')/file_format('tif').
$newQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (f)-->(samp:sample)
MATCH (f)-->(diag:diagnosis)
WHERE s.clinical_study_designation IN ['NCATS-COP01'] and demo.sex in ['Female'] and labels(parent)[0] IN ['diagnosis'] and f.file_type in ['Pathology Report'] and f.file_format IN ['tif']
WITH
        DISTINCT f, parent, c, demo, diag, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
WITH
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN
        coalesce(f.file_name, '') AS `File Name`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_type, '') AS `File Type`,
        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(samp.sample_id, '') AS `Sample ID`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed ,
        coalesce(diag.disease_term,'') AS Diagnosis
        order by f.file_name asc
        limit 100
'@

$ws.Range("B4").Value = $newQuery

# Editing the cell value can trigger an automatic row-height recalculation
# because the cell is styled with wrap text. Restore the original row height
# so only the cell content (and selection, below) change.
$ws.Rows.Item(4).RowHeight = 164.25

# Move the active selection from A5 to C5, as last made by the editor.
$ws.Range("C5").Select()
